$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.7020509200022044
$ws.Range("E2").Value = 0.6591064774027398
$ws.Range("G2").Value = 0.7159537408372849
$ws.Range("I2").Value = 0.6754103730229087
$ws.Range("J2").Value = 0.6736902310684053
$ws.Range("K2").Value = 0.7207704319285049
$ws.Range("L2").Value = 0.6736902310684053
$ws.Range("M2").Value = 0.6812638326675795
$ws.Range("N2").Value = 0.7617707618393961
$ws.Range("O2").Value = 0.7746930423729357
$ws.Range("P2").Value = 0.7617707618393961
$ws.Range("Q2").Value = 0.7611321514088113
$ws.Range("R2").Value = 0.7725692061313202
$ws.Range("S2").Value = 0.7778172554919252
$ws.Range("T2").Value = 0.7725692061313202
$ws.Range("U2").Value = 0.7712238612542054
$ws.Range("V2").Value = 0.7682452528025625
$ws.Range("W2").Value = 0.7826359218639817
$ws.Range("X2").Value = 0.7682452528025625
$ws.Range("Y2").Value = 0.7674291100701333
$ws.Range("B3").Value = 0.8112102493708534
$ws.Range("C3").Value = 0.8207966514864495
$ws.Range("D3").Value = 0.8112102493708534
$ws.Range("E3").Value = 0.8111071707280153
$ws.Range("F3").Value = 0.8112102493708534
$ws.Range("G3").Value = 0.8242138907844732
$ws.Range("H3").Value = 0.8112102493708534
$ws.Range("I3").Value = 0.8125024364594374
$ws.Range("J3").Value = 0.8197208876687258
$ws.Range("K3").Value = 0.8284665768499497
$ws.Range("L3").Value = 0.8197208876687258
$ws.Range("M3").Value = 0.8199223638417941
$ws.Range("N3").Value = 0.8455730954015099
$ws.Range("O3").Value = 0.8545911030801108
$ws.Range("P3").Value = 0.8455730954015099
$ws.Range("Q3").Value = 0.8456102159964825
$ws.Range("R3").Value = 0.8541294898192634
$ws.Range("S3").Value = 0.8619421711243861
$ws.Range("T3").Value = 0.8541294898192634
$ws.Range("U3").Value = 0.8541764935796747
$ws.Range("V3").Value = 0.8541066117593228
$ws.Range("W3").Value = 0.8615440293584173
$ws.Range("X3").Value = 0.8541066117593228
$ws.Range("Y3").Value = 0.8539940367270453
$ws.Range("B4").Value = 0.8368794326241135
$ws.Range("C4").Value = 0.8428092006304425
$ws.Range("D4").Value = 0.8368794326241135
$ws.Range("E4").Value = 0.837042800867278
$ws.Range("F4").Value = 0.8669869595058339
$ws.Range("G4").Value = 0.8733827777519755
$ws.Range("H4").Value = 0.8669869595058339
$ws.Range("I4").Value = 0.8672795407235553
$ws.Range("J4").Value = 0.8433310455273393
$ws.Range("K4").Value = 0.849635434910638
$ws.Range("L4").Value = 0.8433310455273393
$ws.Range("M4").Value = 0.8436890525666965
$ws.Range("N4").Value = 0.8498741706703271
$ws.Range("O4").Value = 0.8566625970168582
$ws.Range("P4").Value = 0.8498741706703271
$ws.Range("Q4").Value = 0.8494232240342564
$ws.Range("R4").Value = 0.8563257835735529
$ws.Range("S4").Value = 0.8625464660246662
$ws.Range("T4").Value = 0.8563257835735529
$ws.Range("U4").Value = 0.8560758551987654
$ws.Range("V4").Value = 0.8563257835735529
$ws.Range("W4").Value = 0.8625464660246662
$ws.Range("X4").Value = 0.8563257835735529
$ws.Range("Y4").Value = 0.8560758551987654
$ws.Range("F5").Value = 0.8476778769160376
$ws.Range("G5").Value = 0.8568907361680995
$ws.Range("H5").Value = 0.8476778769160376
$ws.Range("I5").Value = 0.8461502014655187
$ws.Range("B6").Value = 0.8476549988560971
$ws.Range("C6").Value = 0.8540070616832077
$ws.Range("D6").Value = 0.8476549988560971
$ws.Range("E6").Value = 0.8474796938251774
$ws.Range("F6").Value = 0.8734614504690003
$ws.Range("G6").Value = 0.8782816257816257
$ws.Range("H6").Value = 0.8734614504690003
$ws.Range("I6").Value = 0.8735016000842328
$ws.Range("J6").Value = 0.8305193319606496
$ws.Range("K6").Value = 0.8429270116841419
$ws.Range("L6").Value = 0.8305193319606496
$ws.Range("M6").Value = 0.8326587172502233
$ws.Range("N6").Value = 0.858361930908259
$ws.Range("O6").Value = 0.8631156795488881
$ws.Range("P6").Value = 0.858361930908259
$ws.Range("Q6").Value = 0.8576040212183778
$ws.Range("R6").Value = 0.8648821779913064
$ws.Range("S6").Value = 0.8712407940098277
$ws.Range("T6").Value = 0.8648821779913064
$ws.Range("U6").Value = 0.8643766410690196
$ws.Range("V6").Value = 0.8519103180050331
$ws.Range("W6").Value = 0.8554689458440858
$ws.Range("X6").Value = 0.8519103180050331
$ws.Range("Y6").Value = 0.8509369948730532
